$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 449; this shifts existing rows 449-520 down to 450-521
$ws.Rows.Item(449).Insert()

# Populate the newly inserted row 449 with data (same record as old row 449,
# but with an updated Fecha (D) and Volumen (J))
$ws.Cells.Item(449, 1).Value = 9
$ws.Cells.Item(449, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(449, 3).Value = "Metropolitana"
$ws.Cells.Item(449, 4).Value = 44984
$ws.Cells.Item(449, 5).Value = 13
$ws.Cells.Item(449, 6).Value = 100112052
$ws.Cells.Item(449, 7).Value = "Albahaca"
$ws.Cells.Item(449, 8).Value = "Sin especificar"
$ws.Cells.Item(449, 9).Value = "Primera"
$ws.Cells.Item(449, 10).Value = 250
$ws.Cells.Item(449, 11).Value = 3000
$ws.Cells.Item(449, 12).Value = 3500
$ws.Cells.Item(449, 13).Value = 3250
$ws.Cells.Item(449, 14).Value = "`$/docena de matas"
$ws.Cells.Item(449, 15).Value = "Región Metropolitana"
$ws.Cells.Item(449, 16).Value = 542
$ws.Cells.Item(449, 17).Value = 6
$ws.Cells.Item(449, 18).Value = "Hortaliza"
